$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "2015 - 2023"
    3  = "2015 - 2023"
    4  = "2007 - 2014"
    5  = "2007 - 2014"
    6  = "2015 - 2023"
    7  = "2007 - 2014"
    8  = "1999 - 2006"
    10 = "2015 - 2023"
    11 = "2015 - 2023"
    12 = "2015 - 2023"
    13 = "2015 - 2023"
    14 = "2015 - 2023"
    15 = "2007 - 2014"
    16 = "2007 - 2014"
    17 = "2007 - 2014"
    19 = "2007 - 2014"
    20 = "2007 - 2014"
    21 = "2007 - 2014"
    22 = "2007 - 2014"
    24 = "2007 - 2014"
    25 = "2007 - 2014"
    26 = "2015 - 2023"
    27 = "2015 - 2023"
    28 = "2015 - 2023"
    29 = "2007 - 2014"
    30 = "2007 - 2014"
    31 = "2015 - 2023"
    32 = "2007 - 2014"
    33 = "2015 - 2023"
    34 = "2015 - 2023"
    35 = "2007 - 2014"
    36 = "2007 - 2014"
    37 = "2007 - 2014"
    38 = "1999 - 2006"
    39 = "2007 - 2014"
    40 = "1999 - 2006"
    41 = "1999 - 2006"
    42 = "1999 - 2006"
    44 = "2007 - 2014"
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}
